$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 25, pushing existing rows 25-28 down to 26-29
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the same repeating data as the block,
# but with its own Fecha (date) and Volumen values.
$ws.Cells.Item(25, 1).Value = 8
$ws.Cells.Item(25, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(25, 3).Value = "Coquimbo"
$ws.Cells.Item(25, 4).Value = 45077
$ws.Cells.Item(25, 5).Value = 4
$ws.Cells.Item(25, 6).Value = 100112039
$ws.Cells.Item(25, 7).Value = "Ciboulette"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 760
$ws.Cells.Item(25, 11).Value = 2000
$ws.Cells.Item(25, 12).Value = 2500
$ws.Cells.Item(25, 13).Value = 2250
$ws.Cells.Item(25, 14).Value = "`$/docena de atados"
$ws.Cells.Item(25, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(25, 16).Value = 750
$ws.Cells.Item(25, 17).Value = 3
$ws.Cells.Item(25, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D
$ws.Cells.Item(25, 4).NumberFormat = $ws.Cells.Item(26, 4).NumberFormat
